$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Give the trailing row-7 cell the same "hyperlink" font as column A by routing it
#     through Hyperlinks.Add, then drop its (placeholder) value/link again. ---
$ws.Hyperlinks.Add($ws.Range("A7"), "https://placeholder.invalid")
$ws.Range("A7").WrapText = $true
$ws.Range("A7").ClearContents()

# --- Drop every existing hyperlink (old rows 2 & 3, plus the A7 placeholder above) so
#     the sheet can be rebuilt cleanly with the new movie list. ---
$ws.Range("A2").Hyperlinks.Delete()

# --- New movie rows: url (column A, hyperlinked) + title (column B) ---
$urls = @(
    "https://s3.phim1280.tv/20240531/uwgSpqxG/2000kb/hls/index.m3u8",
    "https://vip.opstream17.com/20240530/9293_d85ee323/3000k/hls/mixed.m3u8",
    "https://vip.opstream17.com/20240530/9425_202c891d/3000k/hls/mixed.m3u8",
    "https://vip.opstream11.com/20240211/51446_7dcc38d0/3000k/hls/mixed.m3u8",
    "https://vip.opstream15.com/20220306/1009_3c3a176a/1000k/hls/mixed.m3u8"
)
$titles = @(
    "Vầng Trăng Máu (Killers of the Flower Moon)",
    "Quý Cô Thích Chiều (Lady Libertine)",
    "Tình yêu, bên kia bức tường (Love, Divided)",
    "Ác Nữ (Lost in Perfection)",
    "Tiền và tình yêu (Money and Love)"
)
$heights = @(30, 45, 45, 45, 45)

for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = 2 + $i
    $aCell = $ws.Cells.Item($row, 1)
    $bCell = $ws.Cells.Item($row, 2)

    $aCell.Value = $urls[$i]
    $bCell.Value = $titles[$i]

    $ws.Hyperlinks.Add($aCell, $urls[$i])
    $aCell.WrapText = $true
    $bCell.WrapText = $true

    $ws.Rows.Item($row).RowHeight = $heights[$i]
}

# --- Selection matches the author's saved cursor position ---
$ws.Range("F7").Select() | Out-Null
